# "Estrutura da pesquisa refeita"
# Turns the raw B3 URL into a real hyperlink run, adds a short
# "Artigo científico" heading paragraph, and appends a new reference
# (ResearchGate URL in Arial) plus a trailing blank paragraph.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------
# 1) Convert the existing "B3 setores" URL paragraph into a hyperlink
# ---------------------------------------------------------------
$b3Url = "https://www.b3.com.br/pt_br/produtos-e-servicos/negociacao/renda-variavel/acoes/consultas/classificacao-setorial/"

$urlPara = $d.Paragraphs(3)
$urlRange = $d.Range($urlPara.Range.Start, $urlPara.Range.End - 1)
$d.Hyperlinks.Add($urlRange, $b3Url, "", "", $b3Url) | Out-Null

# ---------------------------------------------------------------
# 2) New paragraph: "Artigo científico" (justified)
# ---------------------------------------------------------------
$artigoText = "Artigo cient" + [char]0x00ED + "fico"
$insertPoint = $d.Range($d.Content.End, $d.Content.End)
$artigoXml = "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/></w:pPr><w:r><w:t>$artigoText</w:t></w:r></w:p>"
$insertPoint.InsertXML($artigoXml) | Out-Null

# ---------------------------------------------------------------
# 3) New paragraph: ResearchGate URL, Arial font, justified
# ---------------------------------------------------------------
$rgUrl = "https://www.researchgate.net/publication/380091062_Eficiencia_do_Uso_da_Inteligencia_Artificial_no_Desenvolvimento_de_Software"
$insertPoint2 = $d.Range($d.Content.End, $d.Content.End)
$rgXml = "<w:p $wNs><w:pPr><w:jc w:val=`"both`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr><w:t>$rgUrl</w:t></w:r></w:p>"
$insertPoint2.InsertXML($rgXml) | Out-Null

# ---------------------------------------------------------------
# 4) Trailing empty paragraph
# ---------------------------------------------------------------
$insertPoint3 = $d.Range($d.Content.End, $d.Content.End)
$emptyXml = "<w:p $wNs/>"
$insertPoint3.InsertXML($emptyXml) | Out-Null

# ---------------------------------------------------------------
# 5) Register the "Hyperlink" and "Unresolved Mention" character
#    styles (Word normally auto-creates these the first time a
#    hyperlink is inserted into a document).
# ---------------------------------------------------------------
$hyperlinkStyle = $d.Styles.Add("Hyperlink", 2)
$hyperlinkStyle.BaseStyle = $d.Styles("Fontepargpadro")
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.UnhideWhenUsed = $true
$hyperlinkStyle.Font.TextColor.ObjectThemeColor = 10
$hyperlinkStyle.Font.Underline = 1

$mentionStyle = $d.Styles.Add("MenoPendente", 2)
$mentionStyle.NameLocal = "Unresolved Mention"
$mentionStyle.BaseStyle = $d.Styles("Fontepargpadro")
$mentionStyle.Priority = 99
$mentionStyle.UnhideWhenUsed = $true
$mentionStyle.Font.Color = 6053472
